# Apply the "I0 and IF added" change:
# - Add header cells I1="I0", J1="IF" (same style as other header cells)
# - Fill I2:J49 with the corresponding numeric data
# - The sheet dimension automatically expands to A1:J49

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style used by the other header cells (e.g. H1) onto the new headers
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-49
$iValues = @(10,8,8,9,7,6,7,9,7,11,7,7,7,6,7,7,7,8,5,6,6,8,6,6,7,6,7,9,7,7,8,9,8,6,10,8,8,8,9,7,7,6,8,8,7,5,7,5)
$jValues = @(10,8,8,9,7,6,7,9,7,11,7,7,7,6,7,7,7,8,6,6,6,8,6,6,7,7,7,9,7,7,8,9,8,6,10,9,8,8,9,8,7,7,8,8,7,5,7,5)

for ($r = 2; $r -le 49; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]   # Column I
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]  # Column J
}
